$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so numeric-looking
# strings (e.g. "1.00", "36.60", "0.0790") are preserved exactly,
# matching the source inlineStr cells instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.674.18'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.518.53'
$ws.Range("E3").Value = '  -2.08%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '304.44'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").Value = '96.92'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.581'
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("D10").Value = '36.60'
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").Value = '0.0811'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '7.52'
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '2.904.90'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = '2.556.17'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '15.02'
$ws.Range("E16").Value = '  +4.68%  '
$ws.Range("D17").Value = '0.862'
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D18").Value = '42.704.42'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '12.88'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("E20").Value = '  -2.33%  '
$ws.Range("D21").Value = '6.45'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("D22").Value = '71.15'
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("D23").Value = '251.41'
$ws.Range("E23").Value = '  -1.01%  '
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("E25").Value = '  -4.05%  '
$ws.Range("D26").Value = '26.97'
$ws.Range("E26").Value = '  -7.00%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = '2.32'
$ws.Range("E28").Value = '  +10.45%  '
$ws.Range("D29").Value = '10.35'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '38.15'
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").Value = '157.30'
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '3.30'
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0790'
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("E36").Value = '  -4.81%  '
$ws.Range("D37").Value = '18.48'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '24.22'
$ws.Range("E39").Value = '  +5.66%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.119'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = '3.40'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = '2.06'
$ws.Range("E42").Value = '  -2.47%  '
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("D46").Value = '2.034.19'
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("D47").Value = '85.41'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = '8.99'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").Value = '2.759.46'
$ws.Range("E49").Value = '  -2.26%  '
$ws.Range("D50").Value = '0.190'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").Value = '101.60'
$ws.Range("E51").Value = '  -4.59%  '

# Restore default (General) formatting so the cell style matches the
# original workbook (no explicit style index on these cells).
$ws.Range("D2:E51").ClearFormats()

